$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row to the "Items" table for Mushroom
$ws.Range("B13").Value = "Mushroom"
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = "Water"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Wood"
$ws.Range("G13").Value = 2

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item("Items")
$tbl.Resize($ws.Range("B3:G13"))

$ws.Range("H14").Select()
